# Rearranges rows 15-38 (columns A:AY) of the active sheet: the full
# content of each row moves to a different row position, per the mapping
# below (target row -> source row, i.e. "row 15 becomes what row 19 used
# to contain", etc.). Row 14 and all other rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 15
$lastRow  = 38
$lastCol  = "AY"

# target row number => source row number (content currently at the
# source row ends up at the target row)
$mapping = @{
    15 = 19
    16 = 17
    17 = 24
    18 = 15
    19 = 16
    20 = 34
    21 = 26
    22 = 21
    23 = 33
    24 = 31
    25 = 32
    26 = 28
    27 = 22
    28 = 38
    29 = 23
    30 = 29
    31 = 18
    32 = 30
    33 = 25
    34 = 35
    35 = 37
    36 = 20
    37 = 27
    38 = 36
}

# Columns Y/Z/AA/AB hold textual dates/times ("2023-08-24", "00:00").
# Excel's Value2 setter auto-coerces date-shaped strings into date
# serials on write, which would corrupt these text cells. They are
# identical on every affected row anyway, so the two blocks below
# (A:X and AC:AY) simply route around them, leaving them untouched.
$blocks = @(
    @{ Start = "A";  End = "X" },
    @{ Start = "AC"; End = $lastCol }
)

# 1) Snapshot every affected row's values before writing anything back,
#    so source rows read during later iterations are unaffected by
#    earlier writes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($block in $blocks) {
        $rng = $ws.Range($block.Start + $r + ":" + $block.End + $r)
        $rowData[$block.Start] = $rng.Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write each target row using the snapshot of its source row.
for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $mapping[$targetRow]
    $rowData = $snapshot[$sourceRow]
    foreach ($block in $blocks) {
        $destRng = $ws.Range($block.Start + $targetRow + ":" + $block.End + $targetRow)
        $destRng.Value2 = $rowData[$block.Start]
    }
}
